$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.8621994474801085
$ws.Range("C1").Value = 0.3218280216476248
$ws.Range("D1").Value = 0.1378005525198916
$ws.Range("E1").Value = 0.6522780825945866

$ws.Range("B2").Value = 0.8518972959555118
$ws.Range("C2").Value = 0.3646020234935835
$ws.Range("D2").Value = 0.1481027040444882
$ws.Range("E2").Value = 0.6242825976283755

$ws.Range("B3").Value = 0.919355279505717
$ws.Range("C3").Value = 0.1936610613459917
$ws.Range("D3").Value = 0.08064472049428291
$ws.Range("E3").Value = 0.7701979307836657

$ws.Range("B4").Value = 0.9624610274355632
$ws.Range("C4").Value = 0.2338914958918662
$ws.Range("D4").Value = 0.03753897256443682
$ws.Range("E4").Value = 0.7800207965124916

$ws.Range("B5").Value = 0.8785554010948127
$ws.Range("C5").Value = 0.3811338059447518
$ws.Range("D5").Value = 0.1214445989051874
$ws.Range("E5").Value = 0.6361117201774994

$ws.Range("B6").Value = 0.9168950152773334
$ws.Range("C6").Value = 0.3741755108285925
$ws.Range("D6").Value = 0.08310498472266656
$ws.Range("E6").Value = 0.6672328301968279

$ws.Range("B7").Value = 0.8915154350722085
$ws.Range("C7").Value = 0.5419630222041938
$ws.Range("D7").Value = 0.1084845649277914
$ws.Range("E7").Value = 0.5781691403972916

$ws.Range("B8").Value = 0.9506492370615794
$ws.Range("C8").Value = 0.6308630140296889
$ws.Range("D8").Value = 0.04935076293842058
$ws.Range("E8").Value = 0.5829117644360738

$ws.Range("B9").Value = 0.8604519621838554
$ws.Range("C9").Value = 0.3538125892207827
$ws.Range("D9").Value = 0.1395480378161446
$ws.Range("E9").Value = 0.6355768656864891

$ws.Range("B10").Value = 0.7344269705475415
$ws.Range("C10").Value = 0.4308932660169187
$ws.Range("D10").Value = 0.2655730294524585
$ws.Range("E10").Value = 0.5132646773801071

$ws.Range("B11").Value = 0.2664658565376535
$ws.Range("C11").Value = 0.09257409570209964
$ws.Range("D11").Value = 0.7335341434623466
$ws.Range("E11").Value = 0.2438881331580717

$ws.Range("B12").Value = 0.9757411692699818
$ws.Range("C12").Value = 0.5913707020147395
$ws.Range("D12").Value = 0.02425883073001819
$ws.Range("E12").Value = 0.6131451132251298

$ws.Range("B13").Value = 0.9726804557130886
$ws.Range("C13").Value = 0.2539923077867834
$ws.Range("D13").Value = 0.02731954428691143
$ws.Range("E13").Value = 0.7756670034362553

$ws.Range("B14").Value = 0.9542135013741299
$ws.Range("C14").Value = 0.2335038645758613
$ws.Range("D14").Value = 0.04578649862587012
$ws.Range("E14").Value = 0.773579660978391

$ws.Range("B15").Value = 0.9178993337032967
$ws.Range("C15").Value = 0.2116059672908891
$ws.Range("D15").Value = 0.08210066629670333
$ws.Range("E15").Value = 0.7575889839463974

$ws.Range("B16").Value = 0.9273415808237633
$ws.Range("C16").Value = 0.2180218071094901
$ws.Range("D16").Value = 0.07265841917623671
$ws.Range("E16").Value = 0.7613505566246426

$ws.Range("B17").Value = 0.9207039455154958
$ws.Range("C17").Value = 0.1650309777155875
$ws.Range("D17").Value = 0.07929605448450412
$ws.Range("E17").Value = 0.7902828020254258

$ws.Range("B18").Value = 0.8757601387296264
$ws.Range("C18").Value = 0.1598735464828767
$ws.Range("D18").Value = 0.1242398612703735
$ws.Range("E18").Value = 0.7550479458603252

